$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

$ws.Range("B2:E2").ClearContents()
$ws.Range("B3").ClearContents()

$ws.Range("C3").Value = 6.9209747495869554
$ws.Range("D3").Value = 8.4765657727640811
$ws.Range("E3").Value = 3.7829787449636769

$ws.Range("B1:E3").Select() | Out-Null
